$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "28.249.86"
$ws.Cells.Item(2, 5).Value = "  -0.53%  "

$ws.Cells.Item(3, 4).Value = "1.804.67"
$ws.Cells.Item(3, 5).Value = "  -0.73%  "

$ws.Cells.Item(4, 4).Value = "'1.003"
$ws.Cells.Item(4, 5).Value = "  +0.05%  "

$ws.Cells.Item(5, 4).Value = "'314.70"
$ws.Cells.Item(5, 5).Value = "  -0.14%  "

$ws.Cells.Item(7, 4).Value = "'0.5259"
$ws.Cells.Item(7, 5).Value = "  +2.39%  "

$ws.Cells.Item(8, 4).Value = "'0.3821"
$ws.Cells.Item(8, 5).Value = "  -2.62%  "

$ws.Cells.Item(9, 4).Value = "'0.08033"
$ws.Cells.Item(9, 5).Value = "  +1.73%  "

$ws.Cells.Item(10, 2).Value = "OKB"
$ws.Cells.Item(10, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(10, 4).Value = "'41.49"
$ws.Cells.Item(10, 5).Value = "  -0.51%  "

$ws.Cells.Item(11, 2).Value = "Polygon"
$ws.Cells.Item(11, 3).Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Cells.Item(11, 4).Value = "'1.103"
$ws.Cells.Item(11, 5).Value = "  -0.32%  "

$ws.Cells.Item(12, 4).Value = "'6.345"
$ws.Cells.Item(12, 5).Value = "  +1.52%  "

$ws.Cells.Item(13, 4).Value = "'1.003"
$ws.Cells.Item(13, 5).Value = "  +0.09%  "

$ws.Cells.Item(14, 4).Value = "'20.65"
$ws.Cells.Item(14, 5).Value = "  -1.48%  "

$ws.Cells.Item(15, 4).Value = "'7.352"
$ws.Cells.Item(15, 5).Value = "  -1.77%  "

$ws.Cells.Item(16, 4).Value = "1.803.91"

$ws.Cells.Item(17, 4).Value = "'92.61"
$ws.Cells.Item(17, 5).Value = "  +0.08%  "

$ws.Cells.Item(18, 4).Value = "'0.00001098"
$ws.Cells.Item(18, 5).Value = "  -2.41%  "

$ws.Cells.Item(19, 4).Value = "'0.06616"
$ws.Cells.Item(19, 5).Value = "  +0.02%  "

$ws.Cells.Item(20, 5).Value = "  +0.11%  "

$ws.Cells.Item(21, 5).Value = "  -1.45%  "

$ws.Cells.Item(22, 4).Value = "'5.972"
$ws.Cells.Item(22, 5).Value = "  -1.80%  "

$ws.Cells.Item(23, 4).Value = "28.302.88"

$ws.Cells.Item(24, 4).Value = "'11.18"
$ws.Cells.Item(24, 5).Value = "  -0.65%  "

$ws.Cells.Item(25, 4).Value = "'2.238"
$ws.Cells.Item(25, 5).Value = "  -1.29%  "

$ws.Cells.Item(26, 4).Value = "'161.25"
$ws.Cells.Item(26, 5).Value = "  +3.95%  "

$ws.Cells.Item(27, 4).Value = "'20.48"
$ws.Cells.Item(27, 5).Value = "  -2.77%  "

$ws.Cells.Item(28, 4).Value = "2.010.63"
$ws.Cells.Item(28, 5).Value = "  -0.87%  "

$ws.Cells.Item(29, 4).Value = "'2.375"
$ws.Cells.Item(29, 5).Value = "  -1.16%  "

$ws.Cells.Item(30, 4).Value = "'123.36"
$ws.Cells.Item(30, 5).Value = "  -1.73%  "

$ws.Cells.Item(31, 4).Value = "'0.1088"
$ws.Cells.Item(31, 5).Value = "  -1.00%  "

$ws.Cells.Item(32, 4).Value = "'1.060"
$ws.Cells.Item(32, 5).Value = "  -3.81%  "

$ws.Cells.Item(33, 4).Value = "'3.682"
$ws.Cells.Item(33, 5).Value = "  +0.90%  "

$ws.Cells.Item(34, 4).Value = "'5.566"
$ws.Cells.Item(34, 5).Value = "  -1.77%  "

$ws.Cells.Item(35, 4).Value = "'0.07273"
$ws.Cells.Item(35, 5).Value = "  +3.38%  "

$ws.Cells.Item(36, 4).Value = "'12.41"
$ws.Cells.Item(36, 5).Value = "  +10.12%  "

$ws.Cells.Item(37, 2).Value = "VeChain"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(37, 4).Value = "'0.02327"
$ws.Cells.Item(37, 5).Value = "  +0.17%  "

$ws.Cells.Item(38, 4).Value = "'0.2165"
$ws.Cells.Item(38, 5).Value = "  -2.40%  "

$ws.Cells.Item(39, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(39, 4).Value = "'5.160"
$ws.Cells.Item(39, 5).Value = "  -0.31%  "

$ws.Cells.Item(40, 4).Value = "'8.736"
$ws.Cells.Item(40, 5).Value = "  -0.87%  "

$ws.Cells.Item(41, 4).Value = "'0.6234"
$ws.Cells.Item(41, 5).Value = "  -0.06%  "

$ws.Cells.Item(42, 4).Value = "'1.168"
$ws.Cells.Item(42, 5).Value = "  -0.68%  "

$ws.Cells.Item(43, 4).Value = "'1.373"
$ws.Cells.Item(43, 5).Value = "  -1.99%  "

$ws.Cells.Item(44, 2).Value = "Decentraland"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(44, 4).Value = "'0.6068"
$ws.Cells.Item(44, 5).Value = "  +2.98%  "

$ws.Cells.Item(45, 2).Value = "EnergySwap"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(45, 4).Value = "'13.30"
$ws.Cells.Item(45, 5).Value = "  -1.33%  "

$ws.Cells.Item(46, 4).Value = "'3.770"
$ws.Cells.Item(46, 5).Value = "  +0.89%  "

$ws.Cells.Item(47, 4).Value = "'127.13"
$ws.Cells.Item(47, 5).Value = "  +1.97%  "

$ws.Cells.Item(48, 5).Value = "  +2.41%  "

$ws.Cells.Item(49, 4).Value = "'1.934"
$ws.Cells.Item(49, 5).Value = "  -1.77%  "

$ws.Cells.Item(50, 4).Value = "'0.06822"
$ws.Cells.Item(50, 5).Value = "  -0.89%  "

$ws.Cells.Item(51, 4).Value = "'73.13"
$ws.Cells.Item(51, 5).Value = "  -1.75%  "
